# Auto-generated edit script applying scheduled market-price refresh updates
# to the Halicarnassus_Profits workbook sheets (per-sheet H:N market columns).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 799.25
$ws.Range("I32").Value = 799
$ws.Range("J32").Value = 799.5
$ws.Range("K32").Value = 799
$ws.Range("L32").Value = 799.5
$ws.Range("M32").Value = -473
$ws.Range("N32").Value = -1451.5
$ws.Range("H40").Value = 6156.4375
$ws.Range("I40").Value = 5093
$ws.Range("J40").Value = 7928.8335
$ws.Range("K40").Value = 5093
$ws.Range("L40").Value = 7928.8335
$ws.Range("M40").Value = -4918
$ws.Range("N40").Value = -8278.833500000001
$ws.Range("H41").Value = 118
$ws.Range("I41").Value = 118
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 118
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 322
$ws.Range("H80").Value = 360.0909
$ws.Range("I80").Value = 280.83334
$ws.Range("J80").Value = 455.2
$ws.Range("K80").Value = 842.5000200000001
$ws.Range("L80").Value = 1365.6
$ws.Range("M80").Value = 155.4999799999999
$ws.Range("N80").Value = -3361.6
$ws.Range("H83").Value = 360.0909
$ws.Range("I83").Value = 280.83334
$ws.Range("J83").Value = 455.2
$ws.Range("K83").Value = 2527.50006
$ws.Range("L83").Value = 4096.8
$ws.Range("M83").Value = 2464.49994
$ws.Range("N83").Value = -14080.8
$ws.Range("H116").Value = 3536.2
$ws.Range("I116").Value = 3296
$ws.Range("J116").Value = 4096.6665
$ws.Range("K116").Value = 3296
$ws.Range("L116").Value = 4096.6665
$ws.Range("M116").Value = 146
$ws.Range("N116").Value = -10980.6665
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H135").Value = 771
$ws.Range("I135").Value = 873.6667
$ws.Range("J135").Value = 155
$ws.Range("K135").Value = 7863.0003
$ws.Range("L135").Value = 1395
$ws.Range("M135").Value = -5328.0003
$ws.Range("N135").Value = -6465
$ws.Range("H138").Value = 5653.154
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 5653.154
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 16959.462
$ws.Range("N138").Value = -27239.462

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1503.4
$ws.Range("I26").Value = 616
$ws.Range("J26").Value = 2834.5
$ws.Range("K26").Value = 616
$ws.Range("L26").Value = 2834.5
$ws.Range("M26").Value = -286
$ws.Range("N26").Value = -3494.5
$ws.Range("H63").Value = 2219.4
$ws.Range("I63").Value = 2024.5
$ws.Range("J63").Value = 2999
$ws.Range("K63").Value = 2024.5
$ws.Range("L63").Value = 2999
$ws.Range("M63").Value = -1338.5
$ws.Range("N63").Value = -4371
$ws.Range("H66").Value = 2219.4
$ws.Range("I66").Value = 2024.5
$ws.Range("J66").Value = 2999
$ws.Range("K66").Value = 10122.5
$ws.Range("L66").Value = 14995
$ws.Range("M66").Value = -6690.5
$ws.Range("N66").Value = -21859

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5283995
$ws.Range("I105").Value = 8534130
$ws.Range("J105").Value = 2525
$ws.Range("K105").Value = 8534130
$ws.Range("L105").Value = 2525
$ws.Range("M105").Value = -8532383
$ws.Range("N105").Value = -6019
$ws.Range("H132").Value = 150000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 150000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 150000
$ws.Range("N132").Value = -160120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2594.8
$ws.Range("I16").Value = 1993.5
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 1993.5
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -1706.5
$ws.Range("N16").Value = -5574
$ws.Range("H31").Value = 5128.275
$ws.Range("I31").Value = 3092.652
$ws.Range("J31").Value = 7882.353
$ws.Range("K31").Value = 3092.652
$ws.Range("L31").Value = 7882.353
$ws.Range("M31").Value = -2797.652
$ws.Range("N31").Value = -8472.352999999999
$ws.Range("H33").Value = 1063
$ws.Range("I33").Value = 654.8889
$ws.Range("J33").Value = 2899.5
$ws.Range("K33").Value = 654.8889
$ws.Range("L33").Value = 2899.5
$ws.Range("M33").Value = -275.8889
$ws.Range("N33").Value = -3657.5
$ws.Range("H34").Value = 5128.275
$ws.Range("I34").Value = 3092.652
$ws.Range("J34").Value = 7882.353
$ws.Range("K34").Value = 3092.652
$ws.Range("L34").Value = 7882.353
$ws.Range("M34").Value = -2890.652
$ws.Range("N34").Value = -8286.352999999999
$ws.Range("H68").Value = 84688.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 84688.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 84688.5
$ws.Range("N68").Value = -86186.5
$ws.Range("H71").Value = 84688.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 84688.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 254065.5
$ws.Range("N71").Value = -261553.5
$ws.Range("H95").Value = 20250
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 20250
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 20250
$ws.Range("N95").Value = -25742
$ws.Range("H113").Value = 2594.8
$ws.Range("I113").Value = 1993.5
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1993.5
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 176.5
$ws.Range("N113").Value = -9340
$ws.Range("H132").Value = 1793.0714
$ws.Range("I132").Value = 1793.0714
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5379.2142
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2849.2142
$ws.Range("H134").Value = 3134.1155
$ws.Range("I134").Value = 2434.4092
$ws.Range("J134").Value = 6982.5
$ws.Range("K134").Value = 7303.2276
$ws.Range("L134").Value = 20947.5
$ws.Range("M134").Value = -4768.2276
$ws.Range("N134").Value = -26017.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 641.6667
$ws.Range("I107").Value = 395.83334
$ws.Range("J107").Value = 723.6111
$ws.Range("K107").Value = 1187.50002
$ws.Range("L107").Value = 2170.8333
$ws.Range("M107").Value = 732.4999800000001
$ws.Range("N107").Value = -6010.8333

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 31333.334
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 31333.334
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 31333.334
$ws.Range("N26").Value = -31893.334
$ws.Range("H50").Value = 31333.334
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 31333.334
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 31333.334
$ws.Range("N50").Value = -32329.334
$ws.Range("H58").Value = 25000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 25000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 25000
$ws.Range("N58").Value = -25554
$ws.Range("H102").Value = 1075
$ws.Range("I102").Value = 1075
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1075
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 547
$ws.Range("H122").Value = 1112.5
$ws.Range("I122").Value = 1112.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3337.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -887.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6563.1763
$ws.Range("I132").Value = 6563.1763
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 19689.5289
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -17159.5289
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1373.5714
$ws.Range("I136").Value = 1373.5714
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4120.7142
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1570.7142

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H136").Value = 1922.3889
$ws.Range("I136").Value = 1523.6957
$ws.Range("J136").Value = 2627.7693
$ws.Range("K136").Value = 4571.0871
$ws.Range("L136").Value = 7883.3079
$ws.Range("M136").Value = -2021.0871
$ws.Range("N136").Value = -12983.3079
